# Scheduled-runner update: refresh currentAveragePrice / LevePrice / LeveProfit
# columns (H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets
# with newly pulled market-board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 3063.8333
$ws.Range("I39").Value = 125
$ws.Range("J39").Value = 6002.6665
$ws.Range("K39").Value = 375
$ws.Range("L39").Value = 18007.9995
$ws.Range("M39").Value = -79
$ws.Range("N39").Value = -18599.9995

$ws.Range("H62").Value = 16673903
$ws.Range("I62").Value = 22228836
$ws.Range("K62").Value = 22228836
$ws.Range("M62").Value = -22228212

$ws.Range("H65").Value = 16673903
$ws.Range("I65").Value = 22228836
$ws.Range("K65").Value = 111144180
$ws.Range("M65").Value = -111141060

$ws.Range("H113").Value = 2386.6667
$ws.Range("J113").Value = 2232.375
$ws.Range("L113").Value = 2232.375
$ws.Range("N113").Value = -8740.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14962.658
$ws.Range("I32").Value = 13441.208
$ws.Range("J32").Value = 17570.857
$ws.Range("K32").Value = 13441.208
$ws.Range("L32").Value = 17570.857
$ws.Range("M32").Value = -13154.208
$ws.Range("N32").Value = -18144.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 18237
$ws.Range("J119").Value = 18237
$ws.Range("L119").Value = 18237
$ws.Range("N119").Value = -27913

$ws.Range("H120").Value = 81127.5
$ws.Range("J120").Value = 81127.5
$ws.Range("L120").Value = 81127.5
$ws.Range("N120").Value = -90803.5

$ws.Range("H123").Value = 90412.5
$ws.Range("J123").Value = 90412.5
$ws.Range("L123").Value = 90412.5
$ws.Range("N123").Value = -100212.5

$ws.Range("H124").Value = 48750
$ws.Range("J124").Value = 48750
$ws.Range("L124").Value = 48750
$ws.Range("N124").Value = -58570

$ws.Range("H125").Value = 48765
$ws.Range("J125").Value = 48765
$ws.Range("L125").Value = 48765
$ws.Range("N125").Value = -58605

$ws.Range("H134").Value = 8075.5557
$ws.Range("I134").Value = 7022.5
$ws.Range("K134").Value = 21067.5
$ws.Range("M134").Value = -18532.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 150.41667
$ws.Range("I7").Value = 138.70589
$ws.Range("J7").Value = 178.85715
$ws.Range("K7").Value = 138.70589
$ws.Range("L7").Value = 178.85715
$ws.Range("M7").Value = -25.70589000000001
$ws.Range("N7").Value = -404.85715

$ws.Range("H20").Value = 55000
$ws.Range("J20").Value = 55000
$ws.Range("L20").Value = 55000
$ws.Range("N20").Value = -55472

$ws.Range("H30").Value = 55000
$ws.Range("J30").Value = 55000
$ws.Range("L30").Value = 55000
$ws.Range("N30").Value = -55182

$ws.Range("H105").Value = 1756.3334
$ws.Range("I105").Value = 1520.2632
$ws.Range("K105").Value = 1520.2632
$ws.Range("M105").Value = 226.7367999999999

$ws.Range("H123").Value = 60000
$ws.Range("J123").Value = 60000
$ws.Range("L123").Value = 60000
$ws.Range("N123").Value = -69800

$ws.Range("H125").Value = 86238.28999999999
$ws.Range("J125").Value = 86238.28999999999
$ws.Range("L125").Value = 86238.28999999999
$ws.Range("N125").Value = -91158.28999999999

$ws.Range("H128").Value = 55000
$ws.Range("J128").Value = 55000
$ws.Range("L128").Value = 55000
$ws.Range("N128").Value = -64960

$ws.Range("H130").Value = 67050
$ws.Range("J130").Value = 67050
$ws.Range("L130").Value = 67050
$ws.Range("N130").Value = -77090

$ws.Range("H132").Value = 31794.523
$ws.Range("I132").Value = 5291.778
$ws.Range("J132").Value = 190811
$ws.Range("K132").Value = 15875.334
$ws.Range("L132").Value = 572433
$ws.Range("M132").Value = -13345.334
$ws.Range("N132").Value = -577493

$ws.Range("H133").Value = 89999.336
$ws.Range("J133").Value = 89999.336
$ws.Range("L133").Value = 89999.336
$ws.Range("N133").Value = -95059.336

$ws.Range("H134").Value = 6547.72
$ws.Range("I134").Value = 4349.533
$ws.Range("J134").Value = 9845
$ws.Range("K134").Value = 13048.599
$ws.Range("L134").Value = 29535
$ws.Range("M134").Value = -10513.599
$ws.Range("N134").Value = -34605

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1420.1
$ws.Range("I114").Value = 980.2308
$ws.Range("J114").Value = 2237
$ws.Range("K114").Value = 2940.6924
$ws.Range("L114").Value = 6711
$ws.Range("M114").Value = 313.3076000000001
$ws.Range("N114").Value = -13219

$ws.Range("H117").Value = 2439.3809
$ws.Range("J117").Value = 2617.2104
$ws.Range("L117").Value = 7851.6312
$ws.Range("N117").Value = -14735.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H102").Value = 1111.2
$ws.Range("I102").Value = 957.0741
$ws.Range("K102").Value = 957.0741
$ws.Range("M102").Value = 664.9259

$ws.Range("H132").Value = 6487.95
$ws.Range("I132").Value = 4056.6
$ws.Range("J132").Value = 13782
$ws.Range("K132").Value = 12169.8
$ws.Range("L132").Value = 41346
$ws.Range("M132").Value = -9639.799999999999
$ws.Range("N132").Value = -46406

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 79356.84
$ws.Range("I7").Value = 113432.22
$ws.Range("J7").Value = 2687.25
$ws.Range("K7").Value = 113432.22
$ws.Range("L7").Value = 2687.25
$ws.Range("M7").Value = -113320.22
$ws.Range("N7").Value = -2911.25

$ws.Range("H22").Value = 3472.7632
$ws.Range("I22").Value = 2207.5
$ws.Range("J22").Value = 4878.6113
$ws.Range("K22").Value = 2207.5
$ws.Range("L22").Value = 4878.6113
$ws.Range("M22").Value = -1912.5
$ws.Range("N22").Value = -5468.6113

$ws.Range("H27").Value = 3472.7632
$ws.Range("I27").Value = 2207.5
$ws.Range("J27").Value = 4878.6113
$ws.Range("K27").Value = 2207.5
$ws.Range("L27").Value = 4878.6113
$ws.Range("M27").Value = -2100.5
$ws.Range("N27").Value = -5092.6113

$ws.Range("H59").Value = 25000
$ws.Range("J59").Value = 25000
$ws.Range("L59").Value = 25000
$ws.Range("N59").Value = -26308

$ws.Range("H74").Value = 48000
$ws.Range("I74").Value = 48000
$ws.Range("K74").Value = 48000
$ws.Range("M74").Value = -47002

$ws.Range("H77").Value = 48000
$ws.Range("I77").Value = 48000
$ws.Range("K77").Value = 144000
$ws.Range("M77").Value = -139008

$ws.Range("H124").Value = 98463.75
$ws.Range("J124").Value = 98463.75
$ws.Range("L124").Value = 98463.75
$ws.Range("N124").Value = -108283.75

$ws.Range("H125").Value = 177811.62
$ws.Range("J125").Value = 177811.62
$ws.Range("L125").Value = 177811.62
$ws.Range("N125").Value = -187651.62

$ws.Range("H126").Value = 79356.84
$ws.Range("I126").Value = 113432.22
$ws.Range("J126").Value = 2687.25
$ws.Range("K126").Value = 340296.66
$ws.Range("L126").Value = 8061.75
$ws.Range("M126").Value = -337826.66
$ws.Range("N126").Value = -13001.75

$ws.Range("H136").Value = 3035709
$ws.Range("I136").Value = 5132634.5
$ws.Range("J136").Value = 6816.6665
$ws.Range("K136").Value = 15397903.5
$ws.Range("L136").Value = 20449.9995
$ws.Range("M136").Value = -15395353.5
$ws.Range("N136").Value = -25549.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 26342
$ws.Range("J97").Value = 26342
$ws.Range("L97").Value = 26342
$ws.Range("N97").Value = -28324

$ws.Range("H136").Value = 3761059
$ws.Range("I136").Value = 4330289
$ws.Range("J136").Value = 4139.8
$ws.Range("K136").Value = 12990867
$ws.Range("L136").Value = 12419.4
$ws.Range("M136").Value = -12988317
$ws.Range("N136").Value = -17519.4

Write-Output "Leve profit sheets refreshed."
